$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Pan de Azúcar" (Primera
# quality Apio) ahead of the existing history, so insert a fresh row at 512
# and push the rest of the table down by one (row 534 -> 535).
$ws.Rows.Item(512).Insert()

$ws.Cells.Item(512, 1).Value = 3
$ws.Cells.Item(512, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(512, 3).Value = "Coquimbo"
$ws.Cells.Item(512, 4).Value = 44939
$ws.Cells.Item(512, 5).Value = 5
$ws.Cells.Item(512, 6).Value = 100112017
$ws.Cells.Item(512, 7).Value = "Apio"
$ws.Cells.Item(512, 8).Value = "Americana (o)"
$ws.Cells.Item(512, 9).Value = "Primera"
$ws.Cells.Item(512, 10).Value = 130
$ws.Cells.Item(512, 11).Value = 11000
$ws.Cells.Item(512, 12).Value = 11000
$ws.Cells.Item(512, 13).Value = 11000
$ws.Cells.Item(512, 14).Value = "$/docena de matas"
$ws.Cells.Item(512, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(512, 16).Value = 1833
$ws.Cells.Item(512, 17).Value = 6
$ws.Cells.Item(512, 18).Value = "Hortaliza"
